# Update the LR-pairs sheet with the new TPM-derived NATMI output.
# The "ECs" sending cluster is now included alongside Inflammatory-Mac and
# Resolving-Mac, and all previously-present rows have refreshed expression
# values, so every data row (2-16) is rewritten in full.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("ECs", "Il1a", "Il1rap", "ECs", 1, 0.3333333333333333, 0.038687, 0.116061, 0.006675966692470672, 0.006675966692470672, 3, 1, 2.175751666666667, 6.527255, 0.1287951194969566, 0.1287951194969566, 0.08417330472833333, 0.757559742555, 0.0008598319279144625, 0.0008598319279144623),
    @("ECs", "Il1a", "Il1rap", "FAPs", 1, 0.3333333333333333, 0.038687, 0.116061, 0.006675966692470672, 0.006675966692470672, 3, 1, 5.856929333333333, 17.570788, 0.3467049686454247, 0.3467049686454247, 0.2265870251186667, 2.039283226068, 0.002314590822790944, 0.002314590822790944),
    @("ECs", "Il1a", "Il1rap", "Inflammatory-Mac", 1, 0.3333333333333333, 0.038687, 0.116061, 0.006675966692470672, 0.006675966692470672, 3, 1, 5.405749, 16.217247, 0.3199970378477111, 0.3199970378477111, 0.209132211563, 1.882189904067, 0.002136289566360596, 0.002136289566360596),
    @("ECs", "Il1a", "Il1rap", "MuSCs", 1, 0.3333333333333333, 0.038687, 0.116061, 0.006675966692470672, 0.006675966692470672, 3, 1, 1.577118333333333, 4.731355, 0.0933586067355302, 0.09335860673553019, 0.06101397696166666, 0.5491257926549999, 0.0006232589490218678, 0.0006232589490218677),
    @("ECs", "Il1a", "Il1rap", "Resolving-Mac", 1, 0.3333333333333333, 0.038687, 0.116061, 0.006675966692470672, 0.006675966692470672, 3, 1, 1.877573666666666, 5.632720999999999, 0.1111442672743775, 0.1111442672743775, 0.07263769244233333, 0.6537392319809999, 0.0007419954263828024, 0.0007419954263828022),
    @("Inflammatory-Mac", "Il1a", "Il1rap", "ECs", 3, 1, 2.929373, 8.788119, 0.5055030521317986, 0.5055030521317986, 3, 1, 2.175751666666667, 6.527255, 0.1287951194969566, 0.1287951194969566, 6.373588187038334, 57.362293683345, 0.06510632600539129, 0.06510632600539128),
    @("Inflammatory-Mac", "Il1a", "Il1rap", "FAPs", 3, 1, 2.929373, 8.788119, 0.5055030521317986, 0.5055030521317986, 3, 1, 5.856929333333333, 17.570788, 0.3467049686454247, 0.3467049686454247, 17.15713065197467, 154.414175867772, 0.1752604198395217, 0.1752604198395217),
    @("Inflammatory-Mac", "Il1a", "Il1rap", "Inflammatory-Mac", 3, 1, 2.929373, 8.788119, 0.5055030521317986, 0.5055030521317986, 3, 1, 5.405749, 16.217247, 0.3199970378477111, 0.3199970378477111, 15.835455165377, 142.519096488393, 0.1617594793051526, 0.1617594793051526),
    @("Inflammatory-Mac", "Il1a", "Il1rap", "MuSCs", 3, 1, 2.929373, 8.788119, 0.5055030521317986, 0.5055030521317986, 3, 1, 1.577118333333333, 4.731355, 0.0933586067355302, 0.09335860673553019, 4.619967863471667, 41.579710771245, 0.04719306064758281, 0.0471930606475828),
    @("Inflammatory-Mac", "Il1a", "Il1rap", "Resolving-Mac", 3, 1, 2.929373, 8.788119, 0.5055030521317986, 0.5055030521317986, 3, 1, 1.877573666666666, 5.632720999999999, 0.1111442672743775, 0.1111442672743775, 5.500113604644333, 49.50102244179899, 0.05618376633415021, 0.0561837663341502),
    @("Resolving-Mac", "Il1a", "Il1rap", "ECs", 3, 1, 2.826906, 8.480718, 0.4878209811757308, 0.4878209811757308, 3, 1, 2.175751666666667, 6.527255, 0.1287951194969566, 0.1287951194969566, 6.15064544101, 55.35580896909, 0.06282896156365088, 0.06282896156365086),
    @("Resolving-Mac", "Il1a", "Il1rap", "FAPs", 3, 1, 2.826906, 8.480718, 0.4878209811757308, 0.4878209811757308, 3, 1, 5.856929333333333, 17.570788, 0.3467049686454247, 0.3467049686454247, 16.556988673976, 149.012898065784, 0.1691299579831121, 0.169129957983112),
    @("Resolving-Mac", "Il1a", "Il1rap", "Inflammatory-Mac", 3, 1, 2.826906, 8.480718, 0.4878209811757308, 0.4878209811757308, 3, 1, 5.405749, 16.217247, 0.3199970378477111, 0.3199970378477111, 15.281544282594, 137.533898543346, 0.1561012689761979, 0.1561012689761979),
    @("Resolving-Mac", "Il1a", "Il1rap", "MuSCs", 3, 1, 2.826906, 8.480718, 0.4878209811757308, 0.4878209811757308, 3, 1, 1.577118333333333, 4.731355, 0.0933586067355302, 0.09335860673553019, 4.458365279209999, 40.12528751289, 0.04554228713892553, 0.04554228713892552),
    @("Resolving-Mac", "Il1a", "Il1rap", "Resolving-Mac", 3, 1, 2.826906, 8.480718, 0.4878209811757308, 0.4878209811757308, 3, 1, 1.877573666666666, 5.632720999999999, 0.1111442672743775, 0.1111442672743775, 5.307724263741999, 47.76951837367799, 0.05421850551384451, 0.05421850551384449)
)

$r = 2
foreach ($row in $rows) {
    $c = 1
    foreach ($val in $row) {
        $ws.Cells.Item($r, $c).Value = $val
        $c++
    }
    $r++
}
